$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 183.33333
$ws.Range("I4").Value = 110.90909
$ws.Range("K4").Value = 110.90909
$ws.Range("M4").Value = 3.090909999999994

$ws.Range("H33").Value = 391.06668
$ws.Range("I33").Value = 363.8
$ws.Range("J33").Value = 527.4
$ws.Range("K33").Value = 363.8
$ws.Range("L33").Value = 527.4
$ws.Range("M33").Value = -134.8
$ws.Range("N33").Value = -985.4

$ws.Range("H51").Value = 2944.3333
$ws.Range("I51").Value = 2499
$ws.Range("K51").Value = 2499
$ws.Range("M51").Value = -2015

$ws.Range("H53").Value = 1990.5454
$ws.Range("I53").Value = 2411.2222
$ws.Range("K53").Value = 2411.2222
$ws.Range("M53").Value = -1774.2222

$ws.Range("H62").Value = 37040836
$ws.Range("I62").Value = 37040836
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 37040836
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -37040212

$ws.Range("H65").Value = 37040836
$ws.Range("I65").Value = 37040836
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 185204180
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -185201060

$ws.Range("H86").Value = 3043.5
$ws.Range("I86").Value = 3325.625
$ws.Range("J86").Value = 2667.3333
$ws.Range("K86").Value = 3325.625
$ws.Range("L86").Value = 2667.3333
$ws.Range("M86").Value = -2202.625
$ws.Range("N86").Value = -4913.3333

$ws.Range("H89").Value = 3043.5
$ws.Range("I89").Value = 3325.625
$ws.Range("J89").Value = 2667.3333
$ws.Range("K89").Value = 16628.125
$ws.Range("L89").Value = 13336.6665
$ws.Range("M89").Value = -11012.125
$ws.Range("N89").Value = -24568.6665

$ws.Range("H116").Value = 3173
$ws.Range("J116").Value = 3994.1667
$ws.Range("L116").Value = 3994.1667
$ws.Range("N116").Value = -10878.1667

$ws.Range("H132").Value = 5652.52
$ws.Range("I132").Value = 4651.628
$ws.Range("J132").Value = 11800.857
$ws.Range("K132").Value = 13954.884
$ws.Range("L132").Value = 35402.571
$ws.Range("M132").Value = -11424.884
$ws.Range("N132").Value = -40462.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 904
$ws.Range("I2").Value = 761.25
$ws.Range("J2").Value = 1475
$ws.Range("K2").Value = 761.25
$ws.Range("L2").Value = 1475
$ws.Range("M2").Value = -648.25
$ws.Range("N2").Value = -1701

$ws.Range("I45").Value = 1028.8125
$ws.Range("J45").Value = 1682
$ws.Range("K45").Value = 1028.8125
$ws.Range("L45").Value = 1682
$ws.Range("M45").Value = -651.8125
$ws.Range("N45").Value = -2436

$ws.Range("H63").Value = 2149.5642
$ws.Range("I63").Value = 2041.375
$ws.Range("J63").Value = 2322.6667
$ws.Range("K63").Value = 2041.375
$ws.Range("L63").Value = 2322.6667
$ws.Range("M63").Value = -1355.375
$ws.Range("N63").Value = -3694.6667

$ws.Range("H66").Value = 2149.5642
$ws.Range("I66").Value = 2041.375
$ws.Range("J66").Value = 2322.6667
$ws.Range("K66").Value = 10206.875
$ws.Range("L66").Value = 11613.3335
$ws.Range("M66").Value = -6774.875
$ws.Range("N66").Value = -18477.3335

$ws.Range("H74").Value = 954.35297
$ws.Range("I74").Value = 481.6
$ws.Range("K74").Value = 481.6
$ws.Range("M74").Value = 392.4

$ws.Range("H77").Value = 954.35297
$ws.Range("I77").Value = 481.6
$ws.Range("K77").Value = 2408
$ws.Range("M77").Value = 1960

$ws.Range("H110").Value = 2315.7778
$ws.Range("I110").Value = 886
$ws.Range("J110").Value = 4103
$ws.Range("K110").Value = 886
$ws.Range("L110").Value = 4103
$ws.Range("M110").Value = 1159
$ws.Range("N110").Value = -8193

$ws.Range("H116").Value = 904
$ws.Range("I116").Value = 761.25
$ws.Range("J116").Value = 1475
$ws.Range("K116").Value = 761.25
$ws.Range("L116").Value = 1475
$ws.Range("M116").Value = 1532.75
$ws.Range("N116").Value = -6063

$ws.Range("H132").Value = 2952.4285
$ws.Range("I132").Value = 2511.524
$ws.Range("K132").Value = 7534.572
$ws.Range("M132").Value = -5004.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 904
$ws.Range("I3").Value = 761.25
$ws.Range("J3").Value = 1475
$ws.Range("K3").Value = 761.25
$ws.Range("L3").Value = 1475
$ws.Range("M3").Value = -647.25
$ws.Range("N3").Value = -1703

$ws.Range("H107").Value = 1564.6666
$ws.Range("I107").Value = 1290
$ws.Range("J107").Value = 1760.8572
$ws.Range("K107").Value = 1290
$ws.Range("L107").Value = 1760.8572
$ws.Range("M107").Value = 630
$ws.Range("N107").Value = -5600.8572

$ws.Range("H134").Value = 2976.102
$ws.Range("I134").Value = 680.3333
$ws.Range("J134").Value = 9333.615
$ws.Range("K134").Value = 2040.9999
$ws.Range("L134").Value = 28000.845
$ws.Range("M134").Value = 494.0001
$ws.Range("N134").Value = -33070.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1255.2565
$ws.Range("I31").Value = 1082.8
$ws.Range("K31").Value = 1082.8
$ws.Range("M31").Value = -787.8

$ws.Range("H34").Value = 1255.2565
$ws.Range("I34").Value = 1082.8
$ws.Range("K34").Value = 1082.8
$ws.Range("M34").Value = -880.8

$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 4000
$ws.Range("K36").Value = 4000
$ws.Range("M36").Value = -3612

$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3840

$ws.Range("H134").Value = 26317978
$ws.Range("I134").Value = 2348.5
$ws.Range("K134").Value = 7045.5
$ws.Range("M134").Value = -4510.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 624.75
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 766.3333
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 2298.9999
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -2768.9999

$ws.Range("H80").Value = 4499.75
$ws.Range("I80").Value = 2999
$ws.Range("K80").Value = 8997
$ws.Range("M80").Value = -8061

$ws.Range("H83").Value = 4499.75
$ws.Range("I83").Value = 2999
$ws.Range("K83").Value = 26991
$ws.Range("M83").Value = -22311

$ws.Range("H87").Value = 1991.3
$ws.Range("I87").Value = 1207
$ws.Range("J87").Value = 2187.375
$ws.Range("K87").Value = 3621
$ws.Range("L87").Value = 6562.125
$ws.Range("M87").Value = -2373
$ws.Range("N87").Value = -9058.125

$ws.Range("H90").Value = 1991.3
$ws.Range("I90").Value = 1207
$ws.Range("J90").Value = 2187.375
$ws.Range("K90").Value = 10863
$ws.Range("L90").Value = 19686.375
$ws.Range("M90").Value = -4623
$ws.Range("N90").Value = -32166.375

$ws.Range("H97").Value = 871
$ws.Range("I97").Value = 657.5
$ws.Range("K97").Value = 1972.5
$ws.Range("M97").Value = -1476.5

$ws.Range("H140").Value = 23737.334
$ws.Range("I140").Value = 54716.26
$ws.Range("J140").Value = 3440.7932
$ws.Range("K140").Value = 164148.78
$ws.Range("L140").Value = 10322.3796
$ws.Range("M140").Value = -158968.78
$ws.Range("N140").Value = -20682.3796

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1950
$ws.Range("I32").Value = 1950
$ws.Range("K32").Value = 1950
$ws.Range("M32").Value = -1633

$ws.Range("H82").Value = 1940.0667
$ws.Range("I82").Value = 1887.6666
$ws.Range("J82").Value = 2149.6667
$ws.Range("K82").Value = 1887.6666
$ws.Range("L82").Value = 2149.6667
$ws.Range("M82").Value = -1526.6666
$ws.Range("N82").Value = -2871.6667

$ws.Range("H85").Value = 1940.0667
$ws.Range("I85").Value = 1887.6666
$ws.Range("J85").Value = 2149.6667
$ws.Range("K85").Value = 1887.6666
$ws.Range("L85").Value = 2149.6667
$ws.Range("M85").Value = -639.6666
$ws.Range("N85").Value = -4645.6667

$ws.Range("H122").Value = 17859042
$ws.Range("I122").Value = 27779500
$ws.Range("J122").Value = 2220.8
$ws.Range("K122").Value = 83338500
$ws.Range("L122").Value = 6662.400000000001
$ws.Range("M122").Value = -83336050
$ws.Range("N122").Value = -11562.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7749.5
$ws.Range("J15").Value = 7749.5
$ws.Range("L15").Value = 7749.5
$ws.Range("N15").Value = -8325.5

$ws.Range("H62").Value = 83341660
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 83341660
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

$ws.Range("H113").Value = 310.56522
$ws.Range("J113").Value = 472.875
$ws.Range("L113").Value = 1418.625
$ws.Range("N113").Value = -5758.625
